$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-09-03 18:30:54"

for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
